{"js": "// Apply the cover-letter template fixes:\n// 1. Remove justified alignment from the date paragraph ({{CURRENT_DATE}}).\n// 2. Remove justified alignment from the \"Dear Hiring Manager,\" paragraph.\n// 3. Rename the {{JOB_POSITION}} placeholder (body paragraph only) to {{JOB_POSITION_p}}.\n// 4. Rename the {{COMPANY_NAME}} placeholder (body paragraphs only, not the title line)\n//    to {{COMPANY_NAME_p}} in all three body occurrences.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1 & 2: strip the justified alignment on the date / greeting paragraphs ---\n// Paragraph 6: \"{{CURRENT_DATE}}\" ; Paragraph 7: \"Dear Hiring Manager,\"\nconst dateParagraph = paragraphs.items[6];\nconst greetingParagraph = paragraphs.items[7];\ndateParagraph.alignment = Word.Alignment.left;\ngreetingParagraph.alignment = Word.Alignment.left;\nawait context.sync();\n\n// --- 3: {{JOB_POSITION}} -> {{JOB_POSITION_p}} in the \"Please accept...\" paragraph only ---\nconst introParagraph = paragraphs.items[9];\nconst jobPositionResults = introParagraph.search(\"{{JOB_POSITION}}\", { matchCase: true });\njobPositionResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < jobPositionResults.items.length; i++) {\n  jobPositionResults.items[i].insertText(\"{{JOB_POSITION_p}}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 4: {{COMPANY_NAME}} -> {{COMPANY_NAME_p}} in the three body paragraphs\n//        (NOT the \"{{JOB_POSITION}} for {{COMPANY_NAME}}\" title line) ---\nconst companyNameParagraphIndexes = [9, 20, 22];\nfor (const idx of companyNameParagraphIndexes) {\n  const paragraph = paragraphs.items[idx];\n  const companyNameResults = paragraph.search(\"{{COMPANY_NAME}}\", { matchCase: true });\n  companyNameResults.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < companyNameResults.items.length; i++) {\n    companyNameResults.items[i].insertText(\"{{COMPANY_NAME_p}}\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the cover-letter template fixes:\n# 1. Remove justified alignment from the date paragraph ({{CURRENT_DATE}}).\n# 2. Remove justified alignment from the \"Dear Hiring Manager,\" paragraph.\n# 3. Rename the {{JOB_POSITION}} placeholder (body paragraph only) to {{JOB_POSITION_p}}.\n# 4. Rename the {{COMPANY_NAME}} placeholder (body paragraphs only, not the title line)\n#    to {{COMPANY_NAME_p}} in all three body occurrences.\n\n$d = $word.ActiveDocument\n\n# --- 1 & 2: strip the justified alignment on the date / greeting paragraphs ---\n# Paragraph 7: \"{{CURRENT_DATE}}\" ; Paragraph 8: \"Dear Hiring Manager,\"\n$dateParagraph = $d.Paragraphs.Item(7)\n$dateParagraph.Alignment = 0\n\n$greetingParagraph = $d.Paragraphs.Item(8)\n$greetingParagraph.Alignment = 0\n\n# --- 3: {{JOB_POSITION}} -> {{JOB_POSITION_p}} in the \"Please accept...\" paragraph only ---\n$introParagraph = $d.Paragraphs.Item(10)\n$introParagraph.Range.Find.Execute(\"{{JOB_POSITION}}\", $false, $false, $false, $false, $false, $true, 0, $false, \"{{JOB_POSITION_p}}\", 2)\n\n# --- 4: {{COMPANY_NAME}} -> {{COMPANY_NAME_p}} in the three body paragraphs\n#        (NOT the \"{{JOB_POSITION}} for {{COMPANY_NAME}}\" title line) ---\n$introParagraph.Range.Find.Execute(\"{{COMPANY_NAME}}\", $false, $false, $false, $false, $false, $true, 0, $false, \"{{COMPANY_NAME_p}}\", 2)\n\n$drawnParagraph = $d.Paragraphs.Item(21)\n$drawnParagraph.Range.Find.Execute(\"{{COMPANY_NAME}}\", $false, $false, $false, $false, $false, $true, 0, $false, \"{{COMPANY_NAME_p}}\", 2)\n\n$welcomeParagraph = $d.Paragraphs.Item(23)\n$welcomeParagraph.Range.Find.Execute(\"{{COMPANY_NAME}}\", $false, $false, $false, $false, $false, $true, 0, $false, \"{{COMPANY_NAME_p}}\", 2)\n"}
